# TrialsSetup update (2025-12-31 16:00)
# The "RECOVERY" trial's Progress value in the Query1 table has moved on
# from 50 to 100.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B11").Value = 100
